$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder rows 3-5 (Dinafex items): rotate so that row3=old row5, row4=old row3, row5=old row4
$ws.Range("D3").Value = "Dinafex 60mg Tablet"
$ws.Range("E3").Value = "30's"
$ws.Range("D4").Value = "Dinafex 120mg Tablet"
$ws.Range("E4").Value = "30's"
$ws.Range("D5").Value = "Dinafex 180mg Tablet"
$ws.Range("E5").Value = "30's"

# Swap rows 8 and 9 (Etorix items)
$ws.Range("D8").Value = "Etorix 120mg Tablet"
$ws.Range("E8").Value = "20's"
$ws.Range("D9").Value = "Etorix 90mg Tablet"
$ws.Range("E9").Value = "30's"

# Swap rows 11 and 12 (Flucloxin items)
$ws.Range("D11").Value = "Flucloxin 500mg Capsule - 36's"
$ws.Range("E11").Value = "36 's"
$ws.Range("D12").Value = "Flucloxin 500mg Capsule"
$ws.Range("E12").Value = "30 's"

# Swap rows 15 and 16 (Ketonic items)
$ws.Range("D15").Value = "Ketonic 30mg IM/IV Injection - 4's"
$ws.Range("E15").Value = "4's"
$ws.Range("D16").Value = "Ketonic 10mg Tablet"
$ws.Range("E16").Value = "20's"

# Swap rows 17 and 18 (Kynol items)
$ws.Range("D17").Value = "Kynol D 25mg Tablet"
$ws.Range("E17").Value = "60 's"
$ws.Range("D18").Value = "Kynol TR 200mg Capsule"
$ws.Range("E18").Value = "30 's"

# Swap rows 26 and 27 (Zithrox items)
$ws.Range("D26").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("E26").Value = "30ml"
$ws.Range("D27").Value = "Zithrox 250mg Tablet - 6's"
$ws.Range("E27").Value = "6's"
